# Auto-generated from the OOXML diff. Applies per-cell text updates
# to Sheet1, matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.727.35'
$ws.Range('E2').Value = '  +1.82%  '
$ws.Range('D3').Value = '1.574.12'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.13'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.493'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.64'
$ws.Range('E8').Value = '  +3.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.23'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0593'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0890'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '1.797.83'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').Value = '1.568.91'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '28.716.68'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.48'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.72'
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.40'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E23').Value = '  -4.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.18'
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.12'
$ws.Range('E25').Value = '  +9.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.97'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.46'
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.105'
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0487'
$ws.Range('E31').Value = '  +2.84%  '
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.22'
$ws.Range('E33').Value = '  -0.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('D35').Value = '1.389.40'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  +2.58%  '
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.37'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.62'
$ws.Range('E39').Value = '  +2.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.527'
$ws.Range('E41').Value = '  -2.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.92'
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.797'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.967'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.36'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('D49').Value = '1.710.73'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.75'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('E51').Value = '  -0.84%  '
